# Swap the full data content (columns B:AD) between pairs of rows.
# Column A (the running id) stays where it is; everything else
# (match id, teams, scores, odds, etc.) moves to the paired row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(46, 47),
    @(135, 136),
    @(139, 140),
    @(178, 179),
    @(183, 184),
    @(190, 191)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("B$($r1):AD$($r1)")
    $rng2 = $ws.Range("B$($r2):AD$($r2)")

    $v1 = $rng1.Value()
    $v2 = $rng2.Value()

    $rng1.Value = $v2
    $rng2.Value = $v1
}
